$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.617.17'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.12%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.643.19'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.64%  '
$ws.Range("E4").Value = '  +0.21%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.75'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.32%  '
$ws.Range("E6").Value = '  +0.84%  '
$ws.Range("E7").Value = '  +0.20%  '
$ws.Range("E8").Value = '  -0.12%  '
$ws.Range("E9").Value = '  +0.75%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.22'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.40%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0843'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.08%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.872.65'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.68%  '
$ws.Range("E13").Value = '  +3.14%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.634.02'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.46%  '
$ws.Range("E15").Value = '  +1.43%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.99'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.36%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.648.71'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.02%  '
$ws.Range("E18").Value = '  +1.40%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '218.43'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.15%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.00'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.23%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.38'
$ws.Range("D21").Style = "Normal"
$ws.Range("E22").Value = '  +1.99%  '
$ws.Range("E23").Value = '  +1.91%  '
$ws.Range("E24").Value = '  +11.46%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.31'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.23%  '
$ws.Range("E26").Value = '  +0.26%  '
$ws.Range("E27").Value = '  -0.47%  '
$ws.Range("E28").Value = '  +3.55%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.86'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.51%  '
$ws.Range("E30").Value = '  +2.70%  '
$ws.Range("E31").Value = '  +0.86%  '
$ws.Range("E32").Value = '  +3.05%  '
$ws.Range("E33").Value = '  +2.54%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.275.36'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +5.56%  '
$ws.Range("E35").Value = '  +2.40%  '
$ws.Range("E36").Value = '  +6.10%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.41'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.30%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.530'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +5.66%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.828'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.56%  '
$ws.Range("E40").Value = '  +0.23%  '
$ws.Range("E41").Value = '  +2.09%  '
$ws.Range("E42").Value = '  -1.51%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.47'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.10%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.783.76'
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '93.15'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.56%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '59.96'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +9.56%  '
$ws.Range("E47").Value = '  +3.11%  '
$ws.Range("E48").Value = '  +0.59%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.82'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.90%  '
$ws.Range("E50").Value = '  +3.83%  '
$ws.Range("E51").Value = '  -0.64%  '
